# Phase-3 budget workbook update:
#  - highlight three "ggeo_*" task rows (B9, B20, B21) with the "Good"
#    (green) cell style on the "begroting" sheet
#  - restore the working selection on "begroting" to B7:B21
#  - add a new "Sheet1" tab (after "begroting") holding the interactor
#    keyboard-shortcut / instruction notes, and make it the active tab

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("begroting")

# Mark the three newly added engine-function rows with the "Good" style,
# matching the rest of the already-highlighted task list in column B.
$ws1.Range("B9").Style  = "Good"
$ws1.Range("B20").Style = "Good"
$ws1.Range("B21").Style = "Good"

# Leave the "begroting" sheet's view selection on the highlighted block.
$ws1.Activate() | Out-Null
$ws1.Range("B7:B21").Select() | Out-Null

# Add the new notes sheet right after "begroting"; Excel names it
# "Sheet1" by default.
$new = $wb.Worksheets.Add($null, $ws1)

$new.Range("A1").Value = "ctrl + x"
$new.Range("B1").Value = "interactor instructions"
$new.Range("A2").Value = "shift + c "
$new.Range("B2").Value = "polygong offset"
$new.Range("A3").Value = "use buffers for nodes and edges"

$new.PageSetup.PaperSize    = 9
$new.PageSetup.Orientation  = 1

# Land the selection/cursor on the next empty row, and keep this new
# sheet as the active tab.
$new.Range("A4").Select() | Out-Null
